$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8733993172645569
$ws.Range("B1").Value = 3.439122438430786
$ws.Range("C1").Value = 2.722492933273315
$ws.Range("D1").Value = 2.491852998733521
$ws.Range("E1").Value = 1.98189640045166
